$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the product name (B1) and description (B2), removing the test-case
# inter-dependency by giving this product its own unique name/description
# rather than reusing values from product 4480's scenario.
$ws.Range("B1").Value = "4480-RBI-SUBMITLOANON02JAN-INDCOLLSHEETON09JAN-1st"
$ws.Range("B2").Value = "448b"

# Move the active selection to B13
$ws.Range("B13").Select()
